{"js": "// Apply the resume text edits described in the commit \"More edits of the resumes!\".\nconst body = context.document.body;\n\n// ---------------------------------------------------------------\n// Helper: replace the first search hit's full text with new text,\n// preserving the existing run formatting (search result ranges are\n// simple Word.Range objects so insertText/replace keeps formatting).\n// ---------------------------------------------------------------\nasync function replaceOnce(findText, replaceText) {\n  const results = body.search(findText, { matchCase: true, matchWildcards: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n  results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Job title tweak.\nawait replaceOnce(\n  \"Linux Sysadmin, Level 3\",\n  \"Linux Systems Administrator, Level 3\"\n);\n\n// 2) \"support tickets for our\" -> \"support tickets opened for our\"\nawait replaceOnce(\n  \"At this job, I am one of the people answering support tickets for our\",\n  \"At this job, I am one of the people answering support tickets opened for our\"\n);\n\n// 3) Fix double comma typo.\nawait replaceOnce(\n  \". I was hired as a Linux neophyte,, but have been able to expand my knowledge\",\n  \". I was hired as a Linux neophyte, but have been able to expand my knowledge\"\n);\n\n// 4) Expand on the \"senior admin\" sentence.\nawait replaceOnce(\n  \"ime to a point where I can do my job more than adequately (as the senior admin on shift) and guide others\",\n  \"ime to a point where I have been the senior administrator on shift for over a year now - I help and guide others, as well as make changes myself\"\n);\n\n// 5) Add \"listed above.\" to the website/github bullet.\nawait replaceOnce(\n  \": these are on my website and github\",\n  \": these are on my website and github, listed above.\"\n);\n\n// 6) \"Assisting colleagues...\" bullet is unchanged in visible text (only run\n// splitting in the source diff, a cosmetic no-op), nothing to do here.\n\n// 7) Expand the TEFL certificate paragraph text.\nawait replaceOnce(\n  \"arned this certificate in order to expand the ability to talk to groups and explore the world. Many of these skills continue to assist me in my current day to day life.\",\n  \"arned this certificate in order to expand my ability to talk to groups and explore the world - I started it being completely unable to give a presentation, and ended up passing the class! This helped me in my pursuit of a number of soft skills - talking to groups and teaching - that I continue honing to this day.\"\n);\n\n// 8) Insert a new, empty \"Heading 2\" paragraph right after the certificate\n// paragraph and before the \"January 2010 - May 2012\" Heading 2 paragraph.\n{\n  body.load(\"paragraphs\");\n  await context.sync();\n  const paras = body.paragraphs;\n  paras.load(\"text,style\");\n  await context.sync();\n\n  let certParaIndex = -1;\n  for (let i = 0; i < paras.items.length; i++) {\n    if (paras.items[i].text.indexOf(\"I earned this certificate\") === 0) {\n      certParaIndex = i;\n      break;\n    }\n  }\n  if (certParaIndex === -1) {\n    throw new Error(\"Could not find the TEFL certificate paragraph.\");\n  }\n  const certPara = paras.items[certParaIndex];\n  const newPara = certPara.insertParagraph(\"\", Word.InsertLocation.after);\n  newPara.style = \"Heading 2\";\n  newPara.spaceBefore = 0;\n  newPara.lineSpacing = 12;\n  await context.sync();\n\n  // The \"January 2010 - May 2012\" Heading 2 paragraph (now two after the\n  // certificate paragraph) picks up the same tightened spacing.\n  body.load(\"paragraphs\");\n  await context.sync();\n  const paras2 = body.paragraphs;\n  paras2.load(\"text,style\");\n  await context.sync();\n  for (let i = 0; i < paras2.items.length; i++) {\n    if (paras2.items[i].text === \"January 2010 - May 2012\") {\n      paras2.items[i].spaceBefore = 0;\n      paras2.items[i].lineSpacing = 12;\n      break;\n    }\n  }\n  await context.sync();\n}\n\n// 9) Education bullet: drop \"though I could have gotten it in any\n// concentration.\" and rewrite the following sentence.\nawait replaceOnce(\n  \" bachelor\\u2019s degree in English with a concentration in British Literature, though I could have gotten it in any concentration. Studied math and computer science as well.\",\n  \" bachelor\\u2019s degree in English with a concentration in British Literature. I also completed a number of math and computer science classes before deciding that English better fit my life goals at the time.\"\n);\n\n// 10) Hobbies: \"playing guitar/piano and singing\" -> \"playing music\"\nawait replaceOnce(\n  \"Outside of work I keep busy by playing guitar/piano and singing, purchasing and configuring various electronics, reading a lot of books, and writing - mainly poetry, but some fiction and \",\n  \"Outside of work I keep busy by playing music, purchasing and configuring various electronics, reading a lot of books, and writing - mainly poetry, but some fiction and \"\n);\n\n// 11) Drop \"For tech-related things, \" before \"I also have a personal server\".\nawait replaceOnce(\n  \" For tech-related things, I also have a p\",\n  \" I also have a p\"\n);\n", "ps1": "# Apply the resume text edits described in the commit \"More edits of the resumes!\".\n$d = $word.ActiveDocument\n\nfunction Replace-Text {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $found = $rng.Find.Execute($FindText, $false, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)\n    if (-not $found) {\n        throw \"Text not found: $FindText\"\n    }\n}\n\n# 1) Job title tweak.\nReplace-Text \"Linux Sysadmin, Level 3\" \"Linux Systems Administrator, Level 3\"\n\n# 2) \"support tickets for our\" -> \"support tickets opened for our\"\nReplace-Text \"At this job, I am one of the people answering support tickets for our\" \"At this job, I am one of the people answering support tickets opened for our\"\n\n# 3) Fix double comma typo.\nReplace-Text \". I was hired as a Linux neophyte,, but have been able to expand my knowledge\" \". I was hired as a Linux neophyte, but have been able to expand my knowledge\"\n\n# 4) Expand on the \"senior admin\" sentence.\nReplace-Text \"ime to a point where I can do my job more than adequately (as the senior admin on shift) and guide others\" \"ime to a point where I have been the senior administrator on shift for over a year now - I help and guide others, as well as make changes myself\"\n\n# 5) Add \"listed above.\" to the website/github bullet.\nReplace-Text \": these are on my website and github\" \": these are on my website and github, listed above.\"\n\n# 6) \"Assisting colleagues...\" bullet is unchanged in visible text (only run\n# splitting in the source diff, a cosmetic no-op), nothing to do here.\n\n# 7) Expand the TEFL certificate paragraph text.\nReplace-Text \"arned this certificate in order to expand the ability to talk to groups and explore the world. Many of these skills continue to assist me in my current day to day life.\" \"arned this certificate in order to expand my ability to talk to groups and explore the world - I started it being completely unable to give a presentation, and ended up passing the class! This helped me in my pursuit of a number of soft skills - talking to groups and teaching - that I continue honing to this day.\"\n\n# 8) Insert a new, empty \"Heading 2\" paragraph right after the certificate\n# paragraph and before the \"January 2010 - May 2012\" Heading 2 paragraph.\n$certIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i++\n    if ($p.Range.Text.StartsWith(\"I earned this certificate\")) {\n        $certIndex = $i\n        break\n    }\n}\nif ($certIndex -eq -1) {\n    throw \"Could not find the TEFL certificate paragraph.\"\n}\n$certRange = $d.Paragraphs.Item($certIndex).Range\n$certRange.Collapse(0)\n$certRange.InsertParagraphAfter()\n\n# Tighten spacing on the new blank Heading 2 paragraph and the following\n# \"January 2010 - May 2012\" Heading 2 paragraph (both now sit right after\n# the certificate paragraph).\n$newPara = $d.Paragraphs.Item($certIndex + 1)\n$datePara = $d.Paragraphs.Item($certIndex + 2)\n$newPara.Range.ParagraphFormat.SpaceBefore = 0\n$newPara.Range.ParagraphFormat.LineSpacing = 12\n$datePara.Range.ParagraphFormat.SpaceBefore = 0\n$datePara.Range.ParagraphFormat.LineSpacing = 12\n\n# 9) Education bullet: drop \"though I could have gotten it in any\n# concentration.\" and rewrite the following sentence.\nReplace-Text \" bachelor\u2019s degree in English with a concentration in British Literature, though I could have gotten it in any concentration. Studied math and computer science as well.\" \" bachelor\u2019s degree in English with a concentration in British Literature. I also completed a number of math and computer science classes before deciding that English better fit my life goals at the time.\"\n\n# 10) Hobbies: \"playing guitar/piano and singing\" -> \"playing music\"\nReplace-Text \"Outside of work I keep busy by playing guitar/piano and singing, purchasing and configuring various electronics, reading a lot of books, and writing - mainly poetry, but some fiction and \" \"Outside of work I keep busy by playing music, purchasing and configuring various electronics, reading a lot of books, and writing - mainly poetry, but some fiction and \"\n\n# 11) Drop \"For tech-related things, \" before \"I also have a personal server\".\nReplace-Text \" For tech-related things, I also have a p\" \" I also have a p\"\n"}
